$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country names (column A) per new ordering ---
$ws.Cells.Item(36, 1).Value = 'Indonesia'
$ws.Cells.Item(37, 1).Value = 'Polonia'
$ws.Cells.Item(38, 1).Value = 'Chile'
$ws.Cells.Item(39, 1).Value = 'Ecuador'
$ws.Cells.Item(40, 1).Value = 'Tailandia'
$ws.Cells.Item(41, 1).Value = 'Islandia'
$ws.Cells.Item(42, 1).Value = 'Singapur'
$ws.Cells.Item(46, 1).Value = 'Filipinas'
$ws.Cells.Item(47, 1).Value = 'Barein'
$ws.Cells.Item(48, 1).Value = 'Egipto'
$ws.Cells.Item(49, 1).Value = 'Estonia'
$ws.Cells.Item(50, 1).Value = 'India'
$ws.Cells.Item(51, 1).Value = 'Hong Kong'
$ws.Cells.Item(52, 1).Value = 'Peru'
$ws.Cells.Item(54, 1).Value = 'Sudafrica'
$ws.Cells.Item(55, 1).Value = 'Irak'
$ws.Cells.Item(56, 1).Value = 'Mexico'
$ws.Cells.Item(59, 1).Value = 'Kuwait'
$ws.Cells.Item(60, 1).Value = 'Croacia'
$ws.Cells.Item(61, 1).Value = 'Armenia'
$ws.Cells.Item(62, 1).Value = 'Colombia'
$ws.Cells.Item(63, 1).Value = 'Argentina'
$ws.Cells.Item(120, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(121, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(123, 1).Value = 'Guam'
$ws.Cells.Item(124, 1).Value = 'Montenegro'
$ws.Cells.Item(126, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(127, 1).Value = 'Mauricio'
$ws.Cells.Item(136, 1).Value = 'Etiopia'
$ws.Cells.Item(137, 1).Value = 'Trinidad yTobago'
$ws.Cells.Item(138, 1).Value = 'Mayotte'
$ws.Cells.Item(140, 1).Value = 'Kenia'
$ws.Cells.Item(141, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(142, 1).Value = 'Tanzania'
$ws.Cells.Item(143, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(144, 1).Value = 'Barbados'
$ws.Cells.Item(150, 1).Value = 'San Bartolome'
$ws.Cells.Item(151, 1).Value = 'Congo'
$ws.Cells.Item(152, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(153, 1).Value = 'El Salvador'
$ws.Cells.Item(154, 1).Value = 'Madagascar'
$ws.Cells.Item(155, 1).Value = 'Namibia'
$ws.Cells.Item(156, 1).Value = 'Liberia'
$ws.Cells.Item(157, 1).Value = 'Curazao'
$ws.Cells.Item(158, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(159, 1).Value = 'Isla de Man'
$ws.Cells.Item(160, 1).Value = 'Bermudas'
$ws.Cells.Item(162, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(163, 1).Value = 'Santa Lucia'
$ws.Cells.Item(164, 1).Value = 'Benin'
$ws.Cells.Item(165, 1).Value = 'Guinea'
$ws.Cells.Item(166, 1).Value = 'Mauritania'
$ws.Cells.Item(167, 1).Value = 'Haiti'
$ws.Cells.Item(168, 1).Value = 'Butan'
$ws.Cells.Item(169, 1).Value = 'Zambia'
$ws.Cells.Item(170, 1).Value = 'Groenlandia'
$ws.Cells.Item(171, 1).Value = 'Angola'
$ws.Cells.Item(172, 1).Value = 'Fiyi'
$ws.Cells.Item(173, 1).Value = 'Sudan'
$ws.Cells.Item(174, 1).Value = 'Santa Sede'
$ws.Cells.Item(175, 1).Value = 'Cabo Verde'
$ws.Cells.Item(176, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(177, 1).Value = 'Gambia'
$ws.Cells.Item(178, 1).Value = 'Montserrat'
$ws.Cells.Item(179, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(180, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(182, 1).Value = 'Zimbabue'
$ws.Cells.Item(183, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(184, 1).Value = 'Niger'
$ws.Cells.Item(185, 1).Value = 'Timor Oriental'
$ws.Cells.Item(186, 1).Value = 'Republica del Chad'
$ws.Cells.Item(187, 1).Value = 'Suazilandia'
$ws.Cells.Item(188, 1).Value = 'San Vicente y las Granadinas'

# --- Update numeric statistics (columns B-H) ---
$ws.Cells.Item(15, 2).Value = 2690
$ws.Cells.Item(15, 3).Value = 41
$ws.Cells.Item(15, 5).Value = 2674
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 7
$ws.Cells.Item(17, 2).Value = 1980
$ws.Cells.Item(17, 3).Value = 21
$ws.Cells.Item(17, 5).Value = 1972
$ws.Cells.Item(23, 5).Value = 1008
$ws.Cells.Item(23, 7).Value = 1
$ws.Cells.Item(23, 8).Value = 7
$ws.Cells.Item(31, 2).Value = 519
$ws.Cells.Item(31, 3).Value = 18
$ws.Cells.Item(31, 5).Value = 503
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 3
$ws.Cells.Item(36, 2).Value = 450
$ws.Cells.Item(36, 3).Value = 81
$ws.Cells.Item(36, 4).Value = 20
$ws.Cells.Item(36, 5).Value = 392
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 6
$ws.Cells.Item(36, 8).Value = 38
$ws.Cells.Item(37, 2).Value = 439
$ws.Cells.Item(37, 3).Value = 14
$ws.Cells.Item(37, 4).Value = 13
$ws.Cells.Item(37, 5).Value = 421
$ws.Cells.Item(37, 6).Value = 3
$ws.Cells.Item(37, 8).Value = 5
$ws.Cells.Item(38, 2).Value = 434
$ws.Cells.Item(38, 4).Value = 6
$ws.Cells.Item(38, 5).Value = 428
$ws.Cells.Item(38, 6).Value = 7
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(39, 2).Value = 426
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 3
$ws.Cells.Item(39, 5).Value = 416
$ws.Cells.Item(39, 6).Value = 2
$ws.Cells.Item(39, 8).Value = 7
$ws.Cells.Item(40, 2).Value = 411
$ws.Cells.Item(40, 3).Value = 89
$ws.Cells.Item(40, 4).Value = 42
$ws.Cells.Item(40, 5).Value = 368
$ws.Cells.Item(40, 8).Value = 1
$ws.Cells.Item(41, 2).Value = 409
$ws.Cells.Item(41, 4).Value = 5
$ws.Cells.Item(41, 5).Value = 404
$ws.Cells.Item(41, 6).Value = 1
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(42, 2).Value = 385
$ws.Cells.Item(42, 4).Value = 131
$ws.Cells.Item(42, 5).Value = 252
$ws.Cells.Item(42, 6).Value = 14
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 2
$ws.Cells.Item(46, 2).Value = 307
$ws.Cells.Item(46, 3).Value = 77
$ws.Cells.Item(46, 4).Value = 13
$ws.Cells.Item(46, 5).Value = 275
$ws.Cells.Item(46, 6).Value = 1
$ws.Cells.Item(46, 7).Value = 1
$ws.Cells.Item(46, 8).Value = 19
$ws.Cells.Item(47, 2).Value = 298
$ws.Cells.Item(47, 4).Value = 125
$ws.Cells.Item(47, 5).Value = 172
$ws.Cells.Item(47, 6).Value = 4
$ws.Cells.Item(47, 8).Value = 1
$ws.Cells.Item(48, 2).Value = 285
$ws.Cells.Item(48, 4).Value = 42
$ws.Cells.Item(48, 5).Value = 235
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 8).Value = 8
$ws.Cells.Item(49, 2).Value = 283
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(49, 5).Value = 282
$ws.Cells.Item(49, 6).Value = 1
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(50, 2).Value = 275
$ws.Cells.Item(50, 3).Value = 26
$ws.Cells.Item(50, 4).Value = 23
$ws.Cells.Item(50, 5).Value = 247
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 8).Value = 5
$ws.Cells.Item(51, 2).Value = 273
$ws.Cells.Item(51, 3).Value = 17
$ws.Cells.Item(51, 4).Value = 98
$ws.Cells.Item(51, 5).Value = 171
$ws.Cells.Item(51, 6).Value = 4
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 4
$ws.Cells.Item(52, 2).Value = 263
$ws.Cells.Item(52, 4).Value = 1
$ws.Cells.Item(52, 5).Value = 258
$ws.Cells.Item(52, 6).Value = 5
$ws.Cells.Item(54, 2).Value = 240
$ws.Cells.Item(54, 3).Value = 38
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 5).Value = 240
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(55, 2).Value = 208
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 49
$ws.Cells.Item(55, 5).Value = 142
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 17
$ws.Cells.Item(56, 2).Value = 203
$ws.Cells.Item(56, 3).Value = 39
$ws.Cells.Item(56, 4).Value = 4
$ws.Cells.Item(56, 5).Value = 197
$ws.Cells.Item(56, 6).Value = 1
$ws.Cells.Item(56, 7).Value = 1
$ws.Cells.Item(56, 8).Value = 2
$ws.Cells.Item(59, 2).Value = 176
$ws.Cells.Item(59, 3).Value = 17
$ws.Cells.Item(59, 4).Value = 27
$ws.Cells.Item(59, 5).Value = 149
$ws.Cells.Item(59, 6).Value = 5
$ws.Cells.Item(60, 2).Value = 168
$ws.Cells.Item(60, 3).Value = 38
$ws.Cells.Item(60, 4).Value = 5
$ws.Cells.Item(60, 5).Value = 162
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 8).Value = 1
$ws.Cells.Item(61, 2).Value = 160
$ws.Cells.Item(61, 3).Value = 24
$ws.Cells.Item(61, 5).Value = 159
$ws.Cells.Item(61, 6).Value = 2
$ws.Cells.Item(62, 3).Value = 13
$ws.Cells.Item(62, 4).Value = 1
$ws.Cells.Item(62, 5).Value = 157
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(63, 2).Value = 158
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 3
$ws.Cells.Item(63, 5).Value = 152
$ws.Cells.Item(63, 8).Value = 3
$ws.Cells.Item(76, 2).Value = 92
$ws.Cells.Item(76, 3).Value = 1
$ws.Cells.Item(76, 5).Value = 75
$ws.Cells.Item(120, 3).Value = 4
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(126, 3).Value = 5
$ws.Cells.Item(126, 4).Value = 1
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(127, 3).Value = 2
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = 1
$ws.Cells.Item(141, 3).Value = 3
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(153, 3).Value = 2
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(156, 3).Value = 1
$ws.Cells.Item(172, 3).Value = 1
$ws.Cells.Item(172, 5).Value = 2
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(173, 2).Value = 2
$ws.Cells.Item(173, 8).Value = 1
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(185, 3).Value = 1

# --- Update last-updated timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 10:16"
